$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.784.54"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.241.88"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.26"
$ws.Range("E5").Value = "  -2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "83.85"
$ws.Range("E7").Value = "  -1.82%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.39"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.10"
$ws.Range("E12").Value = "  -10.14%  "
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "2.584.41"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.35"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.24"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "2.238.56"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").Value = "39.705.89"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.83"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.46"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "229.27"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.86"
$ws.Range("E27").Value = "  +4.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.03"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.21"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.87"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.71"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.89"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0707"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.32"
$ws.Range("E37").Value = "  +7.11%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0978"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.69"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.73"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "1.929.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  -9.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.60"
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.16"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "2.457.33"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.70"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.39"
$ws.Range("E51").Value = "  -2.46%  "
